$wb = $excel.ActiveWorkbook

# --- Overview sheet: refresh handoff status + generate-date, and shrink the
#     now-shorter Status columns back down (Excel auto-fit after the edit) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-05 21:09:18"
$wsOverview.Range("E1:F1").ColumnWidth = 16.3

# --- zh-cn sheet: refresh handoff status + handoff datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-05 21:09:14"
$wsZhCn.Range("C1").ColumnWidth = 16.3

# --- de-de sheet: refresh handoff status + handoff datetime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-05 21:09:18"
$wsDeDe.Range("C1").ColumnWidth = 16.3
